$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily clear rows beyond the intended filter boundary so that
# re-applying AutoFilter snaps to the exact target range (A1:D39),
# matching the contiguous-region detection used by Range.AutoFilter().
$ws.Range("A40:F42").ClearContents()
$ws.AutoFilterMode = $false
$ws.Range("A1:D39").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$D`$39"
    }
}

# Rewrite the full data grid (rows 1-49) with the updated videocue list.
# Row 1
$ws.Range('A1').Value = 'item'
$ws.Range('B1').Value = 'run_number'
$ws.Range('C1').Value = 'class'
$ws.Range('D1').Value = 'drawn_often'
$ws.Range('E1').Value = 'quickdraw dataset'
$ws.Range('F1').Value = 'alt quickname name'

# Row 2
$ws.Range('A2').Value = 'a bear'
$ws.Range('B2').Value = 'cdm_run_v3_or_v4'
$ws.Range('C2').Value = 'animate'
$ws.Range('D2').Value = 'no'
$ws.Range('E2').Value = 'yes'

# Row 3
$ws.Range('A3').Value = 'a sheep'
$ws.Range('B3').Value = 'cdm_run_v3_or_v4'
$ws.Range('C3').Value = 'animate'
$ws.Range('D3').Value = 'no'
$ws.Range('E3').Value = 'yes'

# Row 4
$ws.Range('A4').Value = 'a camel'
$ws.Range('B4').Value = 'cdm_run_v5'
$ws.Range('C4').Value = 'animate'
$ws.Range('D4').Value = 'no'
$ws.Range('E4').Value = 'yes'

# Row 5
$ws.Range('A5').Value = 'a tiger'
$ws.Range('B5').Value = 'cdm_run_v5'
$ws.Range('C5').Value = 'animate'
$ws.Range('D5').Value = 'no'
$ws.Range('E5').Value = 'yes'

# Row 6
$ws.Range('A6').Value = 'a whale'
$ws.Range('B6').Value = 'cdm_run_v5'
$ws.Range('C6').Value = 'animate'
$ws.Range('D6').Value = 'no'
$ws.Range('E6').Value = 'yes'

# Row 7
$ws.Range('A7').Value = 'a cow'
$ws.Range('B7').Value = 'cdm_run_v6'
$ws.Range('C7').Value = 'animate'
$ws.Range('D7').Value = 'no'
$ws.Range('E7').Value = 'yes'

# Row 8
$ws.Range('A8').Value = 'an elephant'
$ws.Range('B8').Value = 'cdm_run_v6'
$ws.Range('C8').Value = 'animate'
$ws.Range('D8').Value = 'no'
$ws.Range('E8').Value = 'yes'

# Row 9
$ws.Range('A9').Value = 'an octopus'
$ws.Range('B9').Value = 'cdm_run_v7'
$ws.Range('C9').Value = 'animate'
$ws.Range('D9').Value = 'no'
$ws.Range('E9').Value = 'yes'

# Row 10
$ws.Range('A10').Value = 'a frog'
$ws.Range('B10').Value = 'cdm_run_v3_or_v4'
$ws.Range('C10').Value = 'animate'
$ws.Range('D10').Value = 'no'
$ws.Range('E10').Value = 'yes'

# Row 11
$ws.Range('A11').Value = 'a bee'
$ws.Range('B11').Value = 'cdm_run_v7'
$ws.Range('C11').Value = 'animate'
$ws.Range('D11').Value = 'no'
$ws.Range('E11').Value = 'no'

# Row 12
$ws.Range('A12').Value = 'a spider'
$ws.Range('B12').Value = 'cdm_run_v7'
$ws.Range('C12').Value = 'animate'
$ws.Range('D12').Value = 'no'
$ws.Range('E12').Value = 'yes'

# Row 13
$ws.Range('A13').Value = 'a piano'
$ws.Range('B13').Value = 'cdm_run_v7'
$ws.Range('C13').Value = 'inanimate'
$ws.Range('D13').Value = 'no'
$ws.Range('E13').Value = 'yes'

# Row 14
$ws.Range('A14').Value = 'a chair'
$ws.Range('B14').Value = 'cdm_run_v3_or_v4'
$ws.Range('C14').Value = 'inanimate'
$ws.Range('D14').Value = 'yes'
$ws.Range('E14').Value = 'yes'

# Row 15
$ws.Range('A15').Value = 'a couch'
$ws.Range('B15').Value = 'cdm_run_v3_or_v4'
$ws.Range('C15').Value = 'inanimate'
$ws.Range('D15').Value = 'no'
$ws.Range('E15').Value = 'yes'

# Row 16
$ws.Range('A16').Value = 'a cactus'
$ws.Range('B16').Value = 'cdm_run_v5'
$ws.Range('C16').Value = 'inanimate'
$ws.Range('D16').Value = 'no'
$ws.Range('E16').Value = 'yes'

# Row 17
$ws.Range('A17').Value = 'a bowl'
$ws.Range('B17').Value = 'cdm_run_v3_or_v4'
$ws.Range('C17').Value = 'inanimate'
$ws.Range('D17').Value = 'yes'
$ws.Range('E17').Value = 'yes'

# Row 18
$ws.Range('A18').Value = 'a key'
$ws.Range('B18').Value = 'cdm_run_v3_or_v4'
$ws.Range('C18').Value = 'inanimate'
$ws.Range('D18').Value = 'no'
$ws.Range('E18').Value = 'yes'

# Row 19
$ws.Range('A19').Value = 'a phone'
$ws.Range('B19').Value = 'cdm_run_v3_or_v4'
$ws.Range('C19').Value = 'inanimate'
$ws.Range('D19').Value = 'no'
$ws.Range('E19').Value = 'yes'

# Row 20
$ws.Range('A20').Value = 'a scissors'
$ws.Range('B20').Value = 'cdm_run_v3_or_v4'
$ws.Range('C20').Value = 'inanimate'
$ws.Range('D20').Value = 'no'
$ws.Range('E20').Value = 'yes'

# Row 21
$ws.Range('A21').Value = 'a bottle'
$ws.Range('B21').Value = 'cdm_run_v5'
$ws.Range('C21').Value = 'inanimate'
$ws.Range('D21').Value = 'no'
$ws.Range('E21').Value = 'yes'

# Row 22
$ws.Range('A22').Value = 'a hat'
$ws.Range('B22').Value = 'cdm_run_v5'
$ws.Range('C22').Value = 'inanimate'
$ws.Range('D22').Value = 'no'
$ws.Range('E22').Value = 'yes'

# Row 23
$ws.Range('A23').Value = 'a lamp'
$ws.Range('B23').Value = 'cdm_run_v5'
$ws.Range('C23').Value = 'inanimate'
$ws.Range('D23').Value = 'no'
$ws.Range('E23').Value = 'no'

# Row 24
$ws.Range('A24').Value = 'a watch'
$ws.Range('B24').Value = 'cdm_run_v5'
$ws.Range('C24').Value = 'inanimate'
$ws.Range('D24').Value = 'no'
$ws.Range('E24').Value = 'yes'

# Row 25
$ws.Range('A25').Value = 'an apple'
$ws.Range('B25').Value = 'cdm_run_v6'
$ws.Range('C25').Value = 'inanimate'
$ws.Range('D25').Value = 'yes'
$ws.Range('E25').Value = 'yes'

# Row 26
$ws.Range('A26').Value = 'an ice cream'
$ws.Range('B26').Value = 'cdm_run_v6'
$ws.Range('C26').Value = 'inanimate'
$ws.Range('D26').Value = 'yes'
$ws.Range('E26').Value = 'yes'
$ws.Range('F26').Value = 'cell phone'

# Row 27
$ws.Range('A27').Value = 'a clock'
$ws.Range('B27').Value = 'cdm_run_v7'
$ws.Range('C27').Value = 'inanimate'
$ws.Range('D27').Value = 'no'
$ws.Range('E27').Value = 'yes'

# Row 28
$ws.Range('A28').Value = 'a mushroom'
$ws.Range('B28').Value = 'cdm_run_v7'
$ws.Range('C28').Value = 'inanimate'
$ws.Range('D28').Value = 'no'
$ws.Range('E28').Value = 'yes'
$ws.Range('F28').Value = 'wine bottle'

# Row 29
$ws.Range('A29').Value = 'a bike'
$ws.Range('B29').Value = 'cdm_run_v3_or_v4'
$ws.Range('C29').Value = 'inanimate'
$ws.Range('D29').Value = 'no'
$ws.Range('E29').Value = 'yes'

# Row 30
$ws.Range('A30').Value = 'an airplane'
$ws.Range('B30').Value = 'cdm_run_v3_or_v4'
$ws.Range('C30').Value = 'inanimate'
$ws.Range('D30').Value = 'no'
$ws.Range('E30').Value = 'yes'
$ws.Range('F30').Value = 'floor lamp'

# Row 31
$ws.Range('A31').Value = 'a bird'
$ws.Range('B31').Value = 'cdm_run_v3_or_v4'
$ws.Range('C31').Value = 'animate'
$ws.Range('D31').Value = 'yes'
$ws.Range('E31').Value = 'yes'
$ws.Range('F31').Value = 'wristwatch'

# Row 32
$ws.Range('A32').Value = 'a cat'
$ws.Range('B32').Value = 'cdm_run_v3_or_v4'
$ws.Range('C32').Value = 'animate'
$ws.Range('D32').Value = 'yes'
$ws.Range('E32').Value = 'yes'

# Row 33
$ws.Range('A33').Value = 'a rabbit'
$ws.Range('B33').Value = 'cdm_run_v3_or_v4'
$ws.Range('C33').Value = 'animate'
$ws.Range('D33').Value = 'yes'
$ws.Range('E33').Value = 'yes'

# Row 34
$ws.Range('A34').Value = 'a face'
$ws.Range('B34').Value = 'cdm_run_v7'
$ws.Range('C34').Value = 'animate'
$ws.Range('D34').Value = 'yes'
$ws.Range('E34').Value = 'yes'

# Row 35
$ws.Range('A35').Value = 'a hand'
$ws.Range('B35').Value = 'cdm_run_v7'
$ws.Range('C35').Value = 'animate'
$ws.Range('D35').Value = 'yes'
$ws.Range('E35').Value = 'yes'

# Row 36
$ws.Range('A36').Value = 'a horse'
$ws.Range('B36').Value = 'cdm_run_v6'
$ws.Range('C36').Value = 'animate'
$ws.Range('D36').Value = 'yes'
$ws.Range('E36').Value = 'yes'

# Row 37
$ws.Range('A37').Value = 'a person'
$ws.Range('B37').Value = 'cdm_run_v3_or_v4'
$ws.Range('C37').Value = 'animate'
$ws.Range('D37').Value = 'yes'
$ws.Range('E37').Value = 'yes'

# Row 38
$ws.Range('A38').Value = 'a dog'
$ws.Range('B38').Value = 'cdm_run_v3_or_v4'
$ws.Range('C38').Value = 'animate'
$ws.Range('D38').Value = 'yes'
$ws.Range('E38').Value = 'yes'

# Row 39
$ws.Range('A39').Value = 'a fish'
$ws.Range('B39').Value = 'cdm_run_v3_or_v4'
$ws.Range('C39').Value = 'animate'
$ws.Range('D39').Value = 'yes'
$ws.Range('E39').Value = 'yes'

# Row 40
$ws.Range('A40').Value = 'a snail'
$ws.Range('B40').Value = 'cdm_run_v6'
$ws.Range('C40').Value = 'animate'
$ws.Range('D40').Value = 'yes'
$ws.Range('E40').Value = 'yes'

# Row 41
$ws.Range('A41').Value = 'a bed'
$ws.Range('B41').Value = 'cdm_run_v6'
$ws.Range('C41').Value = 'inanimate'
$ws.Range('D41').Value = 'yes'
$ws.Range('E41').Value = 'yes'

# Row 42
$ws.Range('A42').Value = 'a house'
$ws.Range('B42').Value = 'cdm_run_v3_or_v4'
$ws.Range('C42').Value = 'inanimate'
$ws.Range('D42').Value = 'yes'
$ws.Range('E42').Value = 'yes'

# Row 43
$ws.Range('A43').Value = 'a tree'
$ws.Range('B43').Value = 'cdm_run_v3_or_v4'
$ws.Range('C43').Value = 'inanimate'
$ws.Range('D43').Value = 'yes'
$ws.Range('E43').Value = 'yes'

# Row 44
$ws.Range('A44').Value = 'a cup'
$ws.Range('B44').Value = 'cdm_run_v3_or_v4'
$ws.Range('C44').Value = 'inanimate'
$ws.Range('D44').Value = 'yes'
$ws.Range('E44').Value = 'yes'

# Row 45
$ws.Range('A45').Value = 'a book'
$ws.Range('B45').Value = 'cdm_run_v6'
$ws.Range('C45').Value = 'inanimate'
$ws.Range('D45').Value = 'yes'
$ws.Range('E45').Value = 'yes'

# Row 46
$ws.Range('A46').Value = 'a TV'
$ws.Range('B46').Value = 'cdm_run_v6'
$ws.Range('C46').Value = 'inanimate'
$ws.Range('D46').Value = 'yes'
$ws.Range('E46').Value = 'yes'

# Row 47
$ws.Range('A47').Value = 'a boat'
$ws.Range('B47').Value = 'cdm_run_v3_or_v4'
$ws.Range('C47').Value = 'inanimate'
$ws.Range('D47').Value = 'yes'
$ws.Range('E47').Value = 'yes'

# Row 48
$ws.Range('A48').Value = 'a car'
$ws.Range('B48').Value = 'cdm_run_v3_or_v4'
$ws.Range('C48').Value = 'inanimate'
$ws.Range('D48').Value = 'yes'
$ws.Range('E48').Value = 'yes'

# Row 49
$ws.Range('A49').Value = 'a train'
$ws.Range('B49').Value = 'cdm_run_v3_or_v4'
$ws.Range('C49').Value = 'inanimate'
$ws.Range('D49').Value = 'yes'
$ws.Range('E49').Value = 'yes'

# Restore the view selection to match the saved state.
$ws.Range("I45").Select()
